$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B53 was stored as inline/text string "2" -> convert to real number 2
$ws.Range("B53").Value = 2

# Append new row 54 with the new annotation data
$ws.Range("A54").Value = "Ying Tang"
# B54 must stay as literal text "4" (not get auto-converted to a number)
$ws.Range("B54").Value = "'4"
$ws.Range("B54").ClearFormats()
$ws.Range("C54").Value = "interested"
$ws.Range("D54").Value = "QSN"
$ws.Range("E54").Value = "MET"
$ws.Range("F54").Value = "a0a400ab-cd67-43a0-98e0-d641a379b0a8"
$ws.Range("G54").Value = "B1QRgziT-_annotated.xlsx"
$ws.Range("H54").Value = "I am also interested to hear more about the semantics of the spectral norm of this object (flattened filterbank), which Ian asked about below."
